$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.868.81"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "2.416.61"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "555.49"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.37%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "160.60"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +2.53%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.ClearFormats()
$ws.Range("E9").Value = "  +6.40%  "
$ws.Range("E10").Value = "  -0.02%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.326"
$c.ClearFormats()
$ws.Range("E11").Value = "  -1.87%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.76"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "67.795.76"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  +1.44%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "22.96"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.47%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "10.29"
$c.ClearFormats()
$ws.Range("E16").Value = "  -3.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "334.79"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.99%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.79"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.28%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  +0.09%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.86"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "66.61"
$c.ClearFormats()
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  -0.60%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "0.0₃0801"
$ws.Range("E25").Value = "  -0.15%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.09"
$c.ClearFormats()
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +0.10%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "420.01"
$c.ClearFormats()
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -0.67%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "160.85"
$c.ClearFormats()
$ws.Range("E31").Value = "  +3.74%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "18.95"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -0.10%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "17.66"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.06%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.24%  "
$ws.Range("E36").Value = "  -2.62%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.26"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  +0.12%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.04"
$c.ClearFormats()
$ws.Range("E39").Value = "  -4.42%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.98"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.32"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.62%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "128.05"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.09%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0712"
$c.ClearFormats()
$ws.Range("E43").Value = "  +0.19%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.474"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.16%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.553"
$c.ClearFormats()
$ws.Range("E45").Value = "  -0.52%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0914"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  -5.54%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "16.43"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.ClearFormats()
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0200"
$ws.Range("E51").Value = "  +2.20%  "
